# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (K instead of Strike#, regen std/mean, calc and write s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 6
$ws.Range("G4").Value = 1
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 3
$ws.Range("G7").Value = 1
